$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 52, shifting existing rows 52-69 down to 53-70
$ws.Rows.Item(52).Insert()

# Copy the date number format used by column D (style index 2 in before.xlsx) from row 51 to the new row 52
$ws.Range("D52").NumberFormat = $ws.Range("D51").NumberFormat

# Populate the new row 52 with data
$ws.Range("A52").Value = 1
$ws.Range("B52").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C52").Value = "Arica y Parinacota"
$ws.Range("D52").Value = 45093
$ws.Range("E52").Value = 15
$ws.Range("F52").Value = "Fruta"
$ws.Range("G52").Value = 100108
$ws.Range("H52").Value = "Tropicales y subtropicales"
$ws.Range("I52").Value = 100108001
$ws.Range("J52").Value = "Guayaba"
$ws.Range("K52").Value = "Sin especificar"
$ws.Range("L52").Value = "Primera"
$ws.Range("M52").Value = 250
$ws.Range("N52").Value = 4000
$ws.Range("O52").Value = 5000
$ws.Range("P52").Value = 4500
$ws.Range("Q52").Value = "$/caja 10 kilos"
$ws.Range("R52").Value = "Región de Arica y Parinacota"
$ws.Range("S52").Value = 450
$ws.Range("T52").Value = 10
